# SP sync upsert: 2025-11-12T17:21:34.8246670Z - MSCA_DF_00 - Course List.xlsx
# Applies the "Amendment" sheet edits: institute correction for MSCA_DF_03,
# combined type-of-change label for MSCA_DF_15 / MSCA_DF_25, a couple of
# "New EC" values filled in, and clearing the now-removed new-course row 20.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Amendment")

# MSCA_DF_03: new institute was a placeholder guess ("Bifrost University?"),
# now resolved to UTW.
$ws.Range("F3").Value = "UTW"

# MSCA_DF_15: now also a course change, not just an institute change.
$ws.Range("G6").Value = "Institute change + course change"

# MSCA_DF_22: New EC filled in as 0 (course dropped).
$ws.Range("D8").Value = 0

# MSCA_DF_25: now also a course change, not just an institute change.
$ws.Range("G9").Value = "Institute change + course change"

# MSCA_DF_30: New EC filled in as 0 (course dropped).
$ws.Range("D11").Value = 0

# New "Ethics in Digital Finance" course row: New EC set to 3.
$ws.Range("D19").Value = 3

# Remove the second placeholder "New course" row (ASE, 2EC) entirely.
$ws.Range("A20:H20").ClearContents()

# Leave the cursor where the author left it after editing.
$ws.Range("F4").Select()
